$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Add the new "RTM" worksheet right after "Test Scenarios" (position
#    it will occupy in the final tab order), but fill it in AFTER the
#    Login-sheet edits below so the shared-string table gets the same
#    ordering as the authoring session that produced the target file.
# ---------------------------------------------------------------------
$testScenarios = $wb.Worksheets.Item("Test Scenarios")
$rtm = $wb.Worksheets.Add($null, $testScenarios)
$rtm.Name = "RTM"

# ---------------------------------------------------------------------
# 2. Login sheet - add the second login test case (TC_LF_002) as row 4.
# ---------------------------------------------------------------------
$login = $wb.Worksheets.Item("Login")

$login.Range("A4").Value = "TC_LF_002"
$login.Range("B4").Value = "TS_001 - Login"
$login.Range("C4").Value = "Verifying logging into the application using valid email and invalid password"
$login.Range("D4").Value = "1. application url https://tutorialsninja.com.demo is open and supported in any browser (Chrome, Firefox, Edge, Safari)" + [char]10 + "2. Login credentials for an existing account is created"
$login.Range("E4").Value = "1. Click on 'My Account' link" + [char]10 + "2. Click on 'Login' option" + [char]10 + "3. Enter valid email in the E-mail Address field" + [char]10 + "4. Enter invalid password in the Password Field" + [char]10 + "5. Click on Login button"
$login.Range("F4").Value = "Email - seleniumpanda@gmail.com" + [char]10 + "Password - Selenium@123456"
$login.Range("G4").Value = "1. User navigates to login page" + [char]10 + "2. System does not allow login with a warning message ""Warning: No match for E-Mail Address and/or Password."""
$login.Range("H4").Value = "1. User navigates to login page" + [char]10 + "2. System does not allow login with a warning message ""Warning: No match for E-Mail Address and/or Password."""
$login.Range("I4").Value = "P0"
$login.Range("A4:I4").WrapText = $true
$login.Range("A4:I4").HorizontalAlignment = -4108
$login.Range("A4:I4").VerticalAlignment = -4108
$login.Rows.Item(4).RowHeight = 180

# ---------------------------------------------------------------------
# 3. RTM sheet content (Req No / Req Desc / TestCase ID / Status).
# ---------------------------------------------------------------------
$rtm.Range("A1").Value = "Req No"
$rtm.Range("B1").Value = "Req Desc"
$rtm.Range("C1").Value = "TestCase ID"
$rtm.Range("D1").Value = "Status"
$rtm.Range("A2").Value = "Req_1"
$rtm.Range("B2").Value = "Login to the application"
$rtm.Range("C2").Value = "TC_LF_001, TC_LF_002"
$rtm.Range("D2").Value = "TC_LF_001 - Pass, TC_LF_002 - Pass"

$rtm.Range("A1:D2").HorizontalAlignment = -4108
$rtm.Range("A1:D2").VerticalAlignment = -4108
$rtm.Range("A1:D2").WrapText = $true
$rtm.Rows.Item(1).RowHeight = 32.4
$rtm.Rows.Item(2).RowHeight = 32.4

# ---------------------------------------------------------------------
# 4. Status markers on the Login sheet: PASS (green) for the existing
#    TC_LF_001 case, FAIL (red) for the new TC_LF_002 case.
# ---------------------------------------------------------------------
$login.Range("J3:J4").Font.Bold = $true

$pass = $login.Range("J3")
$pass.Value = "PASS"
$pass.Font.Color = 5287936

$fail = $login.Range("J4")
$fail.Value = "FAIL"
$fail.Font.Color = 255

# ---------------------------------------------------------------------
# 5. Sheet-view bookkeeping to mirror the target workbook: Login stays
#    the tab-selected sheet, with the new J4 cell selected, and the
#    active tab (RTM at position 1, zero-based "2" counting in Excel's
#    activeTab attribute) becomes the visible one when the file opens.
# ---------------------------------------------------------------------
$login.Range("J4").Select()
$appWin = $excel.ActiveWindow
$appWin.ScrollColumn = 2

$wb.Worksheets.Item("Login").Activate()
